$d = $word.ActiveDocument

# Locate the paragraph that reads "Added a reminder note inside the
# reconfirm password impromptu pop-up." so we can insert a brand-new list
# item directly after it (and before the "Imp/UK Gal..." bullet).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r`a", "`r", "`n") -eq "Added a reminder note inside the reconfirm password impromptu pop-up.") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find target paragraph 'Added a reminder note inside the reconfirm password impromptu pop-up.'"
}

# Duplicate the target paragraph (including its paragraph mark) and insert the
# copy right after it. This preserves the exact paragraph/run formatting
# (ListParagraph style, numId=2 list numbering, en-US language) used by the
# rest of that list.
$srcRange = $target.Range.Duplicate
$insertionPoint = $d.Range($target.Range.End, $target.Range.End)
$insertionPoint.FormattedText = $srcRange.FormattedText

# The newly inserted paragraph now immediately follows the target paragraph.
$newPara = $target.Next()

# Replace its text (keep the trailing paragraph mark) with the new note.
$textRange = $newPara.Range
$null = $textRange.MoveEnd(1, -1)
$textRange.Text = "Changed the color of reminder note to RED. Font size changed to 20."
